$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column Q: pricing_interest_rate_type ---
$ws.Cells.Item(1, 17).Value = "pricing_interest_rate_type"

$qValues = @(1, 2, 3, 4, 5, 6, 7, 900, 1, 3)
for ($i = 0; $i -lt $qValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 17).Value = $qValues[$i]
}

# --- Column R: pricing_init_rate_period (data first, header last) ---
$rValues = @(12, 36, 48, 12.5, 1.5, "K", "B", "c", 3.4, 60)
for ($i = 0; $i -lt $rValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 18).Value = $rValues[$i]
}

$ws.Cells.Item(1, 18).Value = "pricing_init_rate_period"

# --- Row heights / view state (best effort) ---
$ws.Rows.Item(1).RowHeight = 51

$ws.Range("S4").Select()
